# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

$sheet2021 = $wb.Worksheets.Item("2021-Q1")

# --- Create the new "2022-Q1" worksheet, positioned between "2021-Q1" and "总计" ---
$newSheet = $wb.Worksheets.Add($null, $sheet2021)
$newSheet.Name = "2022-Q1"

# NOTE: fetch "总计" only AFTER the new sheet has been inserted - a reference
# captured beforehand does not track the sheet across the insert.
$sheetTotal = $wb.Worksheets.Item("总计")

# Header row (bold, thin border, centered/top aligned - same look as the other sheets)
$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row 2
$dataRange = $newSheet.Range("A2")
$dataRange.Font.Bold = $true
$dataRange.Borders.LineStyle = 1
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4160
$newSheet.Range("A2").Value = 0

$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "006555"
$newSheet.Range("C2").Value = "浦银安盛全球智能科技股票（QDII）"
$newSheet.Range("D2").Value = "3.20"
$newSheet.Range("E2").Value = "85.41"
$newSheet.Range("F2").Value = "2.15"
$newSheet.Range("G2").Value = "0.0688"
# text values are set; drop the helper "@" number format so the cells end up
# unstyled, matching the plain (no s=) data cells used elsewhere in the workbook
$newSheet.Range("B2:G2").ClearFormats()
$newSheet.Range("H2").Value = 10

# --- Update the "总计" worksheet: insert a new summary row for 2022-Q1 above 2021-Q1 ---
$sheetTotal.Rows.Item(2).Insert()
$sheetTotal.Range("B2:D2").ClearFormats()

$sheetTotal.Range("A2").Font.Bold = $true
$sheetTotal.Range("A2").Borders.LineStyle = 1
$sheetTotal.Range("A2").HorizontalAlignment = -4108
$sheetTotal.Range("A2").VerticalAlignment = -4160
$sheetTotal.Range("A2").Value = 0

$sheetTotal.Range("B2").Value = "2022-Q1"
$sheetTotal.Range("C2").Value = 1
$sheetTotal.Range("D2").Value = 0.07

# the A column is a recomputed 0-based row index - bump the older (now 2nd) row
$sheetTotal.Range("A3").Value = 1

# restore the originally active tab (adding a sheet otherwise activates it)
$sheet2021.Activate()

Write-Host "2022-Q1 sheet added and zongji updated"
